$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 450; this shifts the old rows 450..479 down to 451..480
# and keeps their contents + formatting intact.
$ws.Rows.Item(450).Insert()

# Populate the newly inserted row 450 with the new weekly record.
$ws.Cells.Item(450, 1).Value = 11
$ws.Cells.Item(450, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(450, 3).Value = "Bíobío"
$ws.Cells.Item(450, 4).Value = 45223
$ws.Cells.Item(450, 5).Value = 8
$ws.Cells.Item(450, 6).Value = 100112009
$ws.Cells.Item(450, 7).Value = "Acelga"
$ws.Cells.Item(450, 8).Value = "Sin especificar"
$ws.Cells.Item(450, 9).Value = "Primera"
$ws.Cells.Item(450, 10).Value = 230
$ws.Cells.Item(450, 11).Value = 600
$ws.Cells.Item(450, 12).Value = 650
$ws.Cells.Item(450, 13).Value = 633
$ws.Cells.Item(450, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(450, 15).Value = "Región de Ñuble"
$ws.Cells.Item(450, 16).Value = 633
$ws.Cells.Item(450, 17).Value = 1
$ws.Cells.Item(450, 18).Value = "Hortaliza"
